$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 17862010
$ws.Range("J76").Value = 4569.8
$ws.Range("L76").Value = 4569.8
$ws.Range("N76").Value = -5199.8
$ws.Range("H79").Value = 17862010
$ws.Range("J79").Value = 4569.8
$ws.Range("L79").Value = 4569.8
$ws.Range("N79").Value = -6753.8
$ws.Range("H133").Value = 91335.234
$ws.Range("J133").Value = 91335.234
$ws.Range("L133").Value = 91335.234
$ws.Range("N133").Value = -101455.234
$ws.Range("H134").Value = 91797.664
$ws.Range("J134").Value = 91797.664
$ws.Range("L134").Value = 91797.664
$ws.Range("N134").Value = -101937.664
$ws.Range("H136").Value = 78935.42999999999
$ws.Range("J136").Value = 78935.42999999999
$ws.Range("L136").Value = 78935.42999999999
$ws.Range("N136").Value = -89135.42999999999
$ws.Range("H138").Value = 2704.5962
$ws.Range("I138").Value = 2103.2856
$ws.Range("J138").Value = 2798.1333
$ws.Range("K138").Value = 6309.8568
$ws.Range("L138").Value = 8394.3999
$ws.Range("M138").Value = -1169.8568
$ws.Range("N138").Value = -18674.3999
$ws.Range("H139").Value = 98315.5
$ws.Range("J139").Value = 98315.5
$ws.Range("L139").Value = 98315.5
$ws.Range("N139").Value = -108595.5
$ws.Range("H140").Value = 70585.3
$ws.Range("J140").Value = 70585.3
$ws.Range("L140").Value = 70585.3
$ws.Range("N140").Value = -80945.3

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1005.3333
$ws.Range("J2").Value = 1008
$ws.Range("L2").Value = 1008
$ws.Range("N2").Value = -1234
$ws.Range("H32").Value = 1708.85
$ws.Range("I32").Value = 1570.237
$ws.Range("J32").Value = 6190.6665
$ws.Range("K32").Value = 1570.237
$ws.Range("L32").Value = 6190.6665
$ws.Range("M32").Value = -1283.237
$ws.Range("N32").Value = -6764.6665
$ws.Range("H45").Value = 8930368
$ws.Range("I45").Value = 1872.4
$ws.Range("K45").Value = 1872.4
$ws.Range("M45").Value = -1495.4
$ws.Range("H74").Value = 41789.64
$ws.Range("I74").Value = 53902.527
$ws.Range("K74").Value = 53902.527
$ws.Range("M74").Value = -53028.527
$ws.Range("H77").Value = 41789.64
$ws.Range("I77").Value = 53902.527
$ws.Range("K77").Value = 269512.635
$ws.Range("M77").Value = -265144.635
$ws.Range("H116").Value = 1005.3333
$ws.Range("J116").Value = 1008
$ws.Range("L116").Value = 1008
$ws.Range("N116").Value = -5596
$ws.Range("H122").Value = 3272.3635
$ws.Range("I122").Value = 3166.2222
$ws.Range("K122").Value = 9498.6666
$ws.Range("M122").Value = -7048.6666
$ws.Range("H132").Value = 1389.0222
$ws.Range("I132").Value = 1120.625
$ws.Range("K132").Value = 3361.875
$ws.Range("M132").Value = -831.875
$ws.Range("H134").Value = 105566
$ws.Range("J134").Value = 105566
$ws.Range("L134").Value = 105566
$ws.Range("N134").Value = -115706
$ws.Range("H139").Value = 88228.8
$ws.Range("J139").Value = 88228.8
$ws.Range("L139").Value = 88228.8
$ws.Range("N139").Value = -98508.8

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1005.3333
$ws.Range("J3").Value = 1008
$ws.Range("L3").Value = 1008
$ws.Range("N3").Value = -1236
$ws.Range("H99").Value = 2643048.8
$ws.Range("I99").Value = 74198.57000000001
$ws.Range("J99").Value = 20625000
$ws.Range("K99").Value = 74198.57000000001
$ws.Range("L99").Value = 20625000
$ws.Range("M99").Value = -72700.57000000001
$ws.Range("N99").Value = -20627996
$ws.Range("H105").Value = 46350.957
$ws.Range("I105").Value = 57587.61
$ws.Range("K105").Value = 57587.61
$ws.Range("M105").Value = -55840.61
$ws.Range("H134").Value = 3311.12
$ws.Range("I134").Value = 1117.6
$ws.Range("K134").Value = 3352.8
$ws.Range("M134").Value = -817.7999999999997
$ws.Range("H135").Value = 92017.5
$ws.Range("J135").Value = 92017.5
$ws.Range("L135").Value = 92017.5
$ws.Range("N135").Value = -102157.5
$ws.Range("H138").Value = 79416.89999999999
$ws.Range("J138").Value = 79907.664
$ws.Range("L138").Value = 79907.664
$ws.Range("N138").Value = -90187.664
$ws.Range("H140").Value = 58680.055
$ws.Range("J140").Value = 45903.168
$ws.Range("L140").Value = 45903.168
$ws.Range("N140").Value = -56263.168

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2856.4634
$ws.Range("I31").Value = 2059.1765
$ws.Range("K31").Value = 2059.1765
$ws.Range("M31").Value = -1764.1765
$ws.Range("H34").Value = 2856.4634
$ws.Range("I34").Value = 2059.1765
$ws.Range("K34").Value = 2059.1765
$ws.Range("M34").Value = -1857.1765
$ws.Range("H86").Value = 117688.664
$ws.Range("I86").Value = 253724.75
$ws.Range("J86").Value = 8859.799999999999
$ws.Range("K86").Value = 253724.75
$ws.Range("L86").Value = 8859.799999999999
$ws.Range("M86").Value = -252601.75
$ws.Range("N86").Value = -11105.8
$ws.Range("H89").Value = 117688.664
$ws.Range("I89").Value = 253724.75
$ws.Range("J89").Value = 8859.799999999999
$ws.Range("K89").Value = 1268623.75
$ws.Range("L89").Value = 44299
$ws.Range("M89").Value = -1263007.75
$ws.Range("N89").Value = -55531
$ws.Range("H132").Value = 1521.0385
$ws.Range("J132").Value = 2083.4614
$ws.Range("L132").Value = 6250.3842
$ws.Range("N132").Value = -11310.3842
$ws.Range("H138").Value = 93244.91
$ws.Range("J138").Value = 95998.5
$ws.Range("L138").Value = 95998.5
$ws.Range("N138").Value = -106278.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1154.0454
$ws.Range("I98").Value = 493
$ws.Range("J98").Value = 1348.4706
$ws.Range("K98").Value = 1479
$ws.Range("L98").Value = 4045.4118
$ws.Range("M98").Value = 19
$ws.Range("N98").Value = -7041.4118

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 64652.777
$ws.Range("J109").Value = 64652.777
$ws.Range("L109").Value = 64652.777
$ws.Range("N109").Value = -66732.777
$ws.Range("H113").Value = 3447144.2
$ws.Range("J113").Value = 8336458
$ws.Range("L113").Value = 8336458
$ws.Range("N113").Value = -8340798
$ws.Range("H126").Value = 3486.75
$ws.Range("I126").Value = 2299.1428
$ws.Range("J126").Value = 11800
$ws.Range("K126").Value = 6897.428400000001
$ws.Range("L126").Value = 35400
$ws.Range("M126").Value = -4427.428400000001
$ws.Range("N126").Value = -40340
$ws.Range("H132").Value = 4395.4443
$ws.Range("I132").Value = 3613.8462
$ws.Range("K132").Value = 10841.5386
$ws.Range("M132").Value = -8311.5386
$ws.Range("H135").Value = 47273
$ws.Range("J135").Value = 47273
$ws.Range("L135").Value = 47273
$ws.Range("N135").Value = -57413
$ws.Range("H140").Value = 98913.39999999999
$ws.Range("J140").Value = 98919.14
$ws.Range("L140").Value = 98919.14
$ws.Range("N140").Value = -109279.14

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 4219.8975
$ws.Range("I55").Value = 1141.25
$ws.Range("J55").Value = 9145.733
$ws.Range("K55").Value = 1141.25
$ws.Range("L55").Value = 9145.733
$ws.Range("M55").Value = -968.25
$ws.Range("N55").Value = -9491.733
$ws.Range("H100").Value = 12123.346
$ws.Range("I100").Value = 12980.474
$ws.Range("K100").Value = 12980.474
$ws.Range("M100").Value = -12439.474

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1827.75
$ws.Range("I113").Value = 1946
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 5838
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -3668
$ws.Range("N113").Value = -7340
$ws.Range("H132").Value = 1176073.2
$ws.Range("I132").Value = 915.04346
$ws.Range("J132").Value = 3106690.2
$ws.Range("K132").Value = 2745.13038
$ws.Range("L132").Value = 9320070.600000001
$ws.Range("M132").Value = -215.1303800000001
$ws.Range("N132").Value = -9325130.600000001
$ws.Range("H138").Value = 107450
$ws.Range("J138").Value = 107450
$ws.Range("L138").Value = 107450
$ws.Range("N138").Value = -117730
